$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18-22 correspond to Time Step values 16-20.
# Set the "System True State" (column B) and "Component 4 Truth" (column F)
# to 0 for these rows, centering the generated list on 0.
foreach ($row in 18..22) {
    $ws.Range("B$row").Value = 0
    $ws.Range("F$row").Value = 0
}
